# Generate Report for Handback
# - Updates "Status" from "Ready for handoff" to "Handed back: in sync with en-US"
#   on the zh-cn and de-de sheets (Overview sheet picks this up automatically
#   since it shares the same string).
# - Fills in the "Latest Target File" (F) / "Latest Handback File" (G) /
#   "Latest Handback DateTime" (H) columns that were previously blank /
#   placeholder, mirroring the already-handed-off source (.md) and
#   translated (.xlf) files, with real hyperlinks + the hyperlink style used
#   elsewhere in the sheet.

$wb = $excel.ActiveWorkbook

$hyperColor = 15570276  # OLE (BGR) encoding of RGB FF6495ED, matches the workbook's existing HyperLink style

# ---- Overview sheet (status rollup for both locales) ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("H2").Value = "2016-03-19 10:13:33"
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/d6b8870cdb49e3459a347a734ba942a996e34297/e2e/1f7b7765-3729-4be4-932d-3416dfea63f9.md", "", "", "1f7b7765-3729-4be4-932d-3416dfea63f9.md")
$wsZh.Range("F2").Font.Underline = $true
$wsZh.Range("F2").Font.Color = $hyperColor
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a214ffda7a01adf17b1cc780fcdf9ff0c5d0e3e4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/1f7b7765-3729-4be4-932d-3416dfea63f9.917c915df52a086becfecd8c08225f13faea3066.zh-cn.xlf", "", "", "1f7b7765-3729-4be4-932d-3416dfea63f9.917c915df52a086becfecd8c08225f13faea3066.zh-cn.xlf")
$wsZh.Range("G2").Font.Underline = $true
$wsZh.Range("G2").Font.Color = $hyperColor

$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("H3").Value = "2016-03-19 10:13:33"
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/d6b8870cdb49e3459a347a734ba942a996e34297/e2e/4f169d50-105f-45f7-9845-96ebae52b681.md", "", "", "4f169d50-105f-45f7-9845-96ebae52b681.md")
$wsZh.Range("F3").Font.Underline = $true
$wsZh.Range("F3").Font.Color = $hyperColor
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a214ffda7a01adf17b1cc780fcdf9ff0c5d0e3e4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/4f169d50-105f-45f7-9845-96ebae52b681.92e03053796661e7778e2cb0c3867da5ee257678.zh-cn.xlf", "", "", "4f169d50-105f-45f7-9845-96ebae52b681.92e03053796661e7778e2cb0c3867da5ee257678.zh-cn.xlf")
$wsZh.Range("G3").Font.Underline = $true
$wsZh.Range("G3").Font.Color = $hyperColor

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("H2").Value = "2016-03-19 10:13:38"
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/d6b8870cdb49e3459a347a734ba942a996e34297/e2e/1f7b7765-3729-4be4-932d-3416dfea63f9.md", "", "", "1f7b7765-3729-4be4-932d-3416dfea63f9.md")
$wsDe.Range("F2").Font.Underline = $true
$wsDe.Range("F2").Font.Color = $hyperColor
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/273f6592c9d398f7e9469826261b95b10e4edb7a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/1f7b7765-3729-4be4-932d-3416dfea63f9.917c915df52a086becfecd8c08225f13faea3066.de-de.xlf", "", "", "1f7b7765-3729-4be4-932d-3416dfea63f9.917c915df52a086becfecd8c08225f13faea3066.de-de.xlf")
$wsDe.Range("G2").Font.Underline = $true
$wsDe.Range("G2").Font.Color = $hyperColor

$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("H3").Value = "2016-03-19 10:13:38"
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/d6b8870cdb49e3459a347a734ba942a996e34297/e2e/4f169d50-105f-45f7-9845-96ebae52b681.md", "", "", "4f169d50-105f-45f7-9845-96ebae52b681.md")
$wsDe.Range("F3").Font.Underline = $true
$wsDe.Range("F3").Font.Color = $hyperColor
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/273f6592c9d398f7e9469826261b95b10e4edb7a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/4f169d50-105f-45f7-9845-96ebae52b681.92e03053796661e7778e2cb0c3867da5ee257678.de-de.xlf", "", "", "4f169d50-105f-45f7-9845-96ebae52b681.92e03053796661e7778e2cb0c3867da5ee257678.de-de.xlf")
$wsDe.Range("G3").Font.Underline = $true
$wsDe.Range("G3").Font.Color = $hyperColor
